$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "RM 232" (row 26) entirely, shifting rows up.
$ws.Rows.Item(26).Delete()

# After the first deletion, the row that contained "SC 92" has shifted up to row 27.
# Delete that row entirely as well, shifting remaining rows up again.
$ws.Rows.Item(27).Delete()

# Fill in previously-missing values and clear others to reflect the new
# "missing data" pattern for column B (A = label col, B = 2nd data col).
$ws.Range("B26").Value = -20.2   # SC 5
$ws.Range("B27").Value = ""      # SC 101 (still missing)
$ws.Range("B30").Value = -19.7   # SC 120
$ws.Range("B32").Value = ""      # SC 193 (now missing)
